$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Move the explanatory note ("编号：教师为工号，学生为学号") from K4 up to K1,
# as a new, larger italic heading, and blank out its old spot.
$noteText = $ws.Range("K4").Value2
$ws.Range("K1").Value = $noteText

$ws.Range("K1").Font.Italic = $true
$ws.Range("K1").Font.Size = 22

$ws.Range("K4").ClearContents()

# Row 1 grows to fit the new, larger heading text.
$ws.Rows.Item(1).RowHeight = 28.55

# Update the last-used selection to reflect where editing ended up.
$ws.Range("K16").Select()
